$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'49.467.07"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.93%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.623.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.75%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.05%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'112.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.22%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'323.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -1.38%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.527"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -1.07%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.02%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.542"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -3.00%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'39.57"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -2.64%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'19.74"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -4.19%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -1.30%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +1.22%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'7.28"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.09%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'3.038.71"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.43%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'2.624.78"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.65%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.854"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -2.89%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'49.353.32"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.03%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = "'InternetComputer(DFINITY)"
$ws.Range("B19").Style = "Normal"
$ws.Range("C19").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").Value = "'12.86"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -3.71%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value = "'ImmutableX"
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").Value = "'2.93"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -4.24%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'6.67"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -2.49%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.0₃0945"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -1.98%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'268.72"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -5.04%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'68.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -5.70%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.53"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -2.78%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'26.19"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -2.25%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +0.02%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'10.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +2.89%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -0.91%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -4.65%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'34.63"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -5.44%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'49.39"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.77%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'5.48"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.90%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +2.05%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -0.06%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'18.85"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -3.46%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'4.89"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +2.76%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -1.61%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'3.09"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.90%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'128.16"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +3.07%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = "'Stellar"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'0.111"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -1.95%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "'EnergySwap"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'22.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.37%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.0324"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +3.50%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -3.75%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'2.047.72"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.92%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'3.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -4.89%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'2.10"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +4.98%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -5.38%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'8.90"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.89%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = "'MultiversX"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'58.66"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.77%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = "'THORChain"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'5.18"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -4.00%  "
$ws.Range("E51").Style = "Normal"
